# Add the 2024/10/20 column (AP) of data to the "合成確率" sheet,
# mirroring the existing layout/formatting used by the preceding date
# columns (B .. AO).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newCol = 42          # column AP
$headerRow = 1
$firstDataRow = 2
$lastDataRow = 53

# --- 1. Header cell (AP1): literal text "2024/10/20", matching the
#        plain style used by the other date headers in row 1 (style
#        index 1: メイリオ font, no fill). The format (font/fill) is
#        copied from the preceding header cell (AO1) so it reuses the
#        workbook's existing style, then the cell is pre-formatted as
#        Text so Excel doesn't re-interpret the slash-separated string
#        as a real date serial number. ---
$ws.Range("AO1").Copy()
$headerCell = $ws.Cells.Item($headerRow, $newCol)
$headerCell.PasteSpecial(-4122)
$headerCell.NumberFormat = "@"
$headerCell.Value = "2024/10/20"
$excel.CutCopyMode = 0

# --- 2. Data values for the new column, row by row. ---
$values = @{
    2 = 197
    3 = 162
    4 = 148.8
    5 = 170.9
    6 = 155.5
    7 = 178.1
    8 = 142.6
    9 = 247.2
    10 = 178.6
    11 = 173.2
    12 = 101.9
    13 = 305.6
    14 = 151.5
    15 = 170
    16 = 155.7
    17 = 159.2
    18 = 132.3
    19 = 149.2
    20 = 165.9
    21 = 161
    22 = 166.1
    23 = 134.9
    24 = 173.2
    25 = 272.6
    26 = 151
    27 = 147.1
    28 = 163.7
    29 = 163
    30 = 122
    31 = 174.3
    32 = 217.1
    33 = 125
    34 = 123
    35 = 165.2
    36 = 233.2
    37 = 214.8
    38 = 152.9
    39 = 170.1
    40 = 187.7
    41 = 141.6
    42 = 126.7
    43 = 246.5
    44 = 111
    45 = 134.7
    46 = 147.2
    47 = 159.5
    48 = 148.4
    49 = 187.3
    50 = 133.4
    51 = 165.1
    52 = 193.6
    53 = 141.5
}

# Reference cells carrying each of the three fill styles already used
# throughout the sheet (plain / yellow-highlight / blue-highlight),
# picked from an existing, untouched column so copying their format
# reuses the workbook's existing style entries instead of creating new
# ones.
$refNormal = $ws.Range("B2")   # style 1: no fill
$refYellow = $ws.Range("D2")   # style 2: fill < 125
$refBlue   = $ws.Range("N2")   # style 3: 125 <= fill < 140

for ($row = $firstDataRow; $row -le $lastDataRow; $row++) {
    $v = $values[$row]

    if ($v -lt 125.0) {
        $ref = $refYellow
    } elseif ($v -lt 140.0) {
        $ref = $refBlue
    } else {
        $ref = $refNormal
    }

    $ref.Copy()
    $target = $ws.Cells.Item($row, $newCol)
    $target.PasteSpecial(-4122)
    $target.Value = $v
}

$excel.CutCopyMode = 0

# --- 3. Column width for AP, matching the width="12" used by every
#        other data column. (11.1666... compensates for this engine's
#        internal character-width -> stored-width padding so the
#        saved OOXML <col> width comes out to exactly 12.) ---
$ws.Columns($newCol).ColumnWidth = 11.1666666666667
